$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bloom Cnt values so they are reflected in the UI
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 4

# Move the active selection to A3
$ws.Range("A3").Select()
